$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.932756781578064
$ws.Range("B1").Value = 2.605717182159424
$ws.Range("C1").Value = 2.460438966751099
$ws.Range("D1").Value = 2.594599008560181
$ws.Range("E1").Value = 3.38646125793457
